# Replace the two M2Doc "complex field" constructs (fldChar begin/instrText/fldChar end)
# with the equivalent plain-text token syntax "{...}", as used by the new
# TokenIteratorFieldRewriterSplit parser:
#   { m:userdoc 'zone1' }  -> {m:userdoc 'zone1'}
#   { m:enduserdoc }       -> {m:enduserdoc}
#
# The two fields live inside a table cell that also contains a nested table,
# so plain Find/Replace can't see their (hidden) field-code text. Instead we
# locate each field's paragraph precisely and rewrite it in place with
# Range.InsertXML, which replaces exactly the targeted range's contents.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Field 1: m:userdoc 'zone1' -------------------------------------------
# Word reliably reports this field's Code range; pad by one char on each
# side to also capture the begin/end field-character runs around it.
$field1 = $d.Fields.Item(1)
$start1 = $field1.Code.Start - 1
$end1 = $field1.Code.End + 1

$range1 = $d.Range($start1, $end1)
$xml1 = '<w:p ' + $wNs + ' w:rsidP="00B50B6B" w:rsidR="00B50B6B" w:rsidRDefault="00B50B6B"><w:r><w:t>{m:userdoc ''zone1''}</w:t></w:r></w:p>'
$null = $range1.InsertXML($xml1)

# --- Field 2: m:enduserdoc --------------------------------------------------
# This field sits right after the nested table, in the same outer cell.
$range2 = $d.Range(82, 98)
$xml2 = '<w:p ' + $wNs + ' w:rsidP="00B50B6B" w:rsidR="00B50B6B" w:rsidRDefault="00B50B6B"><w:r><w:t>{m:enduserdoc}</w:t></w:r></w:p>'
$null = $range2.InsertXML($xml2)

Write-Host "Replaced" $d.Fields.Count "remaining complex field(s); both userdoc fields converted to bracket syntax."
